# Regenerate the "K" (strikeouts) column (column G) on the active sheet with
# freshly computed values (replacing the previous "Strike#" derived values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value (column G), per the regenerated save_data
$newK = @{
    2  = 2
    3  = 2
    4  = 9
    5  = 8
    6  = 3
    7  = 6
    8  = 4
    9  = 3
    10 = 6
    11 = 9
    12 = 3
    13 = 5
    14 = 6
    15 = 5
    16 = 4
    17 = 8
    18 = 7
    19 = 4
    20 = 7
    21 = 7
    22 = 3
    23 = 9
    24 = 2
    25 = 6
    26 = 6
    27 = 6
    28 = 6
    29 = 4
    30 = 3
    31 = 1
    32 = 2
    33 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
